$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-08-05 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-08-06 Tuesday", 2)

# Update the division problems in the table, cell by cell (row, column)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "31÷7="
$t.Cell(1, 2).Range.Text = "75÷3="
$t.Cell(1, 3).Range.Text = "79÷5="
$t.Cell(1, 4).Range.Text = "49÷8="
$t.Cell(1, 5).Range.Text = "37÷8="

$t.Cell(5, 1).Range.Text = "65÷2="
$t.Cell(5, 2).Range.Text = "89÷8="
$t.Cell(5, 3).Range.Text = "88÷5="
$t.Cell(5, 4).Range.Text = "24÷7="
$t.Cell(5, 5).Range.Text = "87÷9="

$t.Cell(9, 1).Range.Text = "96÷3="
$t.Cell(9, 2).Range.Text = "72÷5="
$t.Cell(9, 3).Range.Text = "87÷8="
$t.Cell(9, 4).Range.Text = "77÷4="
$t.Cell(9, 5).Range.Text = "81÷9="

$t.Cell(13, 1).Range.Text = "15÷9="
$t.Cell(13, 2).Range.Text = "56÷9="
$t.Cell(13, 3).Range.Text = "17÷5="
$t.Cell(13, 4).Range.Text = "22÷4="
$t.Cell(13, 5).Range.Text = "89÷3="

$t.Cell(17, 1).Range.Text = "76÷3="
$t.Cell(17, 2).Range.Text = "77÷2="
$t.Cell(17, 3).Range.Text = "33÷4="
$t.Cell(17, 4).Range.Text = "58÷9="
$t.Cell(17, 5).Range.Text = "68÷4="
